# Regenerate the lattice-multiplication worksheet table with a new set of
# 15 exercises (5 rows x 3 columns), matching the freshly generated output
# at commit c8c62b6. Each cell's run contains 5 lines separated by line
# breaks: "A x B", the column digits, the "----" separator, and the two
# row-digit lines.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# vertical-tab char == a <w:br/> line break inside a Word Range.Text
$nl = [char]11

function Set-Cell($row, $col, $line1, $line2, $line4, $line5) {
    $cell = $t.Cell($row, $col)
    $text = $line1 + $nl + $line2 + $nl + "  ----" + $nl + $line4 + $nl + $line5
    $cell.Range.Text = $text
}

Set-Cell 1 1 "37 x 55" "  5    5" "3|    |" "7|    |"
Set-Cell 1 2 "11 x 14" "  1    4" "1|    |" "1|    |"
Set-Cell 1 3 "92 x 41" "  4    1" "9|    |" "2|    |"

Set-Cell 2 1 "48 x 16" "  1    6" "4|    |" "8|    |"
Set-Cell 2 2 "50 x 48" "  4    8" "5|    |" "0|    |"
Set-Cell 2 3 "69 x 39" "  3    9" "6|    |" "9|    |"

Set-Cell 3 1 "75 x 61" "  6    1" "7|    |" "5|    |"
Set-Cell 3 2 "59 x 99" "  9    9" "5|    |" "9|    |"
Set-Cell 3 3 "15 x 97" "  9    7" "1|    |" "5|    |"

Set-Cell 4 1 "59 x 68" "  6    8" "5|    |" "9|    |"
Set-Cell 4 2 "22 x 43" "  4    3" "2|    |" "2|    |"
Set-Cell 4 3 "75 x 69" "  6    9" "7|    |" "5|    |"

Set-Cell 5 1 "90 x 45" "  4    5" "9|    |" "0|    |"
Set-Cell 5 2 "88 x 21" "  2    1" "8|    |" "8|    |"
Set-Cell 5 3 "99 x 60" "  6    0" "9|    |" "9|    |"

Write-Host "Done updating" $t.Rows.Count "x" $t.Columns.Count "table"
